$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 124 (before the current row 125),
# pushing all existing rows 125..146 down to 127..148.
$ws.Rows.Item(125).Insert()
$ws.Rows.Item(125).Insert()

# New row 125: Camote, Primera, week of 2023-03-20 ($/caja 18 kilos)
$ws.Cells.Item(125, 1).Value = 9
$ws.Cells.Item(125, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(125, 3).Value = "Metropolitana"
$ws.Cells.Item(125, 4).Value = 45005
$ws.Cells.Item(125, 5).Value = 13
$ws.Cells.Item(125, 6).Value = 100114002
$ws.Cells.Item(125, 7).Value = "Camote"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 790
$ws.Cells.Item(125, 11).Value = 17000
$ws.Cells.Item(125, 12).Value = 18000
$ws.Cells.Item(125, 13).Value = 17494
$ws.Cells.Item(125, 14).Value = "`$/caja 18 kilos"
$ws.Cells.Item(125, 15).Value = "Perú"
$ws.Cells.Item(125, 16).Value = 972
$ws.Cells.Item(125, 17).Value = 18
$ws.Cells.Item(125, 18).Value = "Hortaliza"

# New row 126: Camote, Primera, week of 2023-03-20 ($/malla 18 kilos)
$ws.Cells.Item(126, 1).Value = 9
$ws.Cells.Item(126, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(126, 3).Value = "Metropolitana"
$ws.Cells.Item(126, 4).Value = 45005
$ws.Cells.Item(126, 5).Value = 13
$ws.Cells.Item(126, 6).Value = 100114002
$ws.Cells.Item(126, 7).Value = "Camote"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 610
$ws.Cells.Item(126, 11).Value = 15000
$ws.Cells.Item(126, 12).Value = 16000
$ws.Cells.Item(126, 13).Value = 15500
$ws.Cells.Item(126, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(126, 15).Value = "Perú"
$ws.Cells.Item(126, 16).Value = 861
$ws.Cells.Item(126, 17).Value = 18
$ws.Cells.Item(126, 18).Value = "Hortaliza"
